# Split the single run of bibliography text in the "Bibliografia" section
# into multiple runs separated by pairs of manual line breaks (<w:br/><w:br/>),
# one break-pair between each reference entry.

$d = $word.ActiveDocument

# Each pair below is the boundary straddling the end of one reference and the
# start of the next one. We search for the exact concatenated text (no break)
# and replace it with the same text but with two manual line breaks (^l^l)
# inserted at the junction. Wildcards are enabled so ^l is interpreted as a
# line break (w:br) rather than literal text.

$boundaries = @(
    @("tos. [s.l.] AMGH, 2013. 472 p.BARBIERI, J. C. Gestão Ambient",
      "tos. [s.l.] AMGH, 2013. 472 p.^l^lBARBIERI, J. C. Gestão Ambient"),
    @("mentos. Editora Saraiva, 2004.ALLEN, D.T.; SHONNARD, D. R., ",
      "mentos. Editora Saraiva, 2004.^l^lALLEN, D.T.; SHONNARD, D. R., "),
    @("studies, Prentice Hall, 2015. AKKUCUK, U. Handbook of Resear",
      "studies, Prentice Hall, 2015. ^l^lAKKUCUK, U. Handbook of Resear"),
    @("s.l.] IGI Global, 2020. 409 p.BOUCHERY, Y.; CORBETT, C. J.; ",
      "s.l.] IGI Global, 2020. 409 p.^l^lBOUCHERY, Y.; CORBETT, C. J.; "),
    @("Publishing, 2017. v. 4. 130 p.SCHMIDT, M.; GIOVANNUCCI, D.; ",
      "Publishing, 2017. v. 4. 130 p.^l^lSCHMIDT, M.; GIOVANNUCCI, D.; "),
    @("Publishing, 2019. v. 2. 304 p.LAVE, L. B.; HENDRICKSON, C. T",
      "Publishing, 2019. v. 2. 304 p.^l^lLAVE, L. B.; HENDRICKSON, C. T"),
    @("s, Editora John Hopkins, 2006.LEITE, P. R. Logística Reversa",
      "s, Editora John Hopkins, 2006.^l^lLEITE, P. R. Logística Reversa")
)

foreach ($pair in $boundaries) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}
